$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Update the Date field (B2). Force the cell to remain text so the
#    date-looking string "2025-07-11" is not silently converted into a
#    real Excel date serial number.
# -----------------------------------------------------------------
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2025-07-11"
$ws.Range("B2").Style = "Normal"

# -----------------------------------------------------------------
# 2. Update the Inclusion Criteria text (B5): replace the last two
#    bullet points with a single "Randomized control trials" bullet.
# -----------------------------------------------------------------
$nl = [char]10
$inclusion = $nl + "    " + [char]0x2022 + " Studies published in English, peer-reviewed journals" + $nl + "    " + [char]0x2022 + " About leptin and Alzheimer" + [char]0x2019 + "s" + $nl + "    " + [char]0x2022 + " Relevant papers available as full text" + $nl + "    " + [char]0x2022 + " Randomized control trials " + $nl + "    "
$ws.Range("B5").Value = $inclusion

# -----------------------------------------------------------------
# 3. Update the Study Type field (B7).
# -----------------------------------------------------------------
$ws.Range("B7").Value = "Randomized control trials"

# -----------------------------------------------------------------
# 4. Reorder the results table (rows 12-55, columns B-G) so that the
#    papers appear sorted by relevance. Column A (the running index
#    1..44) and columns H/I (empty) are left untouched.
# -----------------------------------------------------------------

# Snapshot the current ("before") content of every row, columns B..G,
# using Value2 (reading through .Value on this host mis-behaves).
$snapshot = @{}
for ($r = 12; $r -le 55; $r++) {
    $row = @()
    for ($c = 2; $c -le 7; $c++) {
        $row += $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $row
}

# Mapping of new ("after") row number -> original ("before") row number.
$map = @{}
$map[12] = 22
$map[13] = 15
$map[14] = 17
$map[15] = 51
$map[16] = 36
$map[17] = 45
$map[18] = 49
$map[19] = 19
$map[20] = 50
$map[21] = 41
$map[22] = 39
$map[23] = 14
$map[24] = 48
$map[25] = 43
$map[26] = 25
$map[27] = 35
$map[28] = 27
$map[29] = 26
$map[30] = 40
$map[31] = 13
$map[32] = 44
$map[33] = 20
$map[34] = 16
$map[35] = 12
$map[36] = 18
$map[37] = 54
$map[38] = 37
$map[39] = 21
$map[40] = 38
$map[41] = 23
$map[42] = 31
$map[43] = 30
$map[44] = 29
$map[45] = 42
$map[46] = 53
$map[47] = 24
$map[48] = 34
$map[49] = 47
$map[50] = 32
$map[51] = 33
$map[52] = 46
$map[53] = 52
$map[54] = 55
$map[55] = 28

# Make sure column C (Paper Year) keeps being stored as text instead of
# being auto-converted to a number when we write the 4-digit year back.
$ws.Range("C12:C55").NumberFormat = "@"

for ($r = 12; $r -le 55; $r++) {
    $srcRow = $map[$r]
    $data = $snapshot[$srcRow]
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 2]
    }
}

$ws.Range("C12:C55").Style = "Normal"
